$d = $word.ActiveDocument

# Insert "Team Members: " before the names in the first paragraph, as its
# own run (matching the author's original two-run layout).
$para = $d.Paragraphs(1)
$r = $para.Range
$r.Collapse(1)
$prefix = "Team Members: "
$r.InsertBefore($prefix)

# Drop a "_GoBack" bookmark right after the inserted prefix. Word only
# ever keeps a single "_GoBack" bookmark, so adding it here also removes
# the stale one that previously sat in the final empty paragraph.
$start = $para.Range.Start + $prefix.Length
$bmRange = $d.Range($start, $start)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Reformat the date from "9/24/2019" to "24 September 2019".
$d.Content.Find.Execute("9/24/2019", $true, $false, $false, $false, $false,
                         $true, 1, $false, "24 September 2019", 2)
